$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '40.947.81'
$ws.Cells.Item(2, 5).Value = '  -1.76%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.169.83'
$ws.Cells.Item(3, 5).Value = '  -2.64%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '246.68'
$ws.Cells.Item(5, 5).Value = '  -2.41%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -2.16%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '66.08'
$ws.Cells.Item(7, 5).Value = '  -7.12%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.05%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.565'
$ws.Cells.Item(9, 5).Value = '  -0.87%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '59.06'
$ws.Cells.Item(10, 5).Value = '  +0.24%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0924'
$ws.Cells.Item(11, 5).Value = '  -4.25%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '35.52'
$ws.Cells.Item(12, 5).Value = '  -16.81%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -1.88%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.86'
$ws.Cells.Item(14, 5).Value = '  -1.71%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.490.82'
$ws.Cells.Item(15, 5).Value = '  -2.46%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.853'
$ws.Cells.Item(16, 5).Value = '  -0.30%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '14.25'
$ws.Cells.Item(17, 5).Value = '  -4.72%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.169.57'
$ws.Cells.Item(18, 5).Value = '  -2.69%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '40.838.82'
$ws.Cells.Item(19, 5).Value = '  -1.82%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0₃0936'
$ws.Cells.Item(20, 5).Value = '  -3.44%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.08'
$ws.Cells.Item(21, 5).Value = '  -2.05%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '71.17'
$ws.Cells.Item(22, 5).Value = '  -2.56%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '228.95'
$ws.Cells.Item(23, 5).Value = '  -2.53%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.12'
$ws.Cells.Item(24, 5).Value = '  -6.99%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'Dai'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.00'
$ws.Cells.Item(25, 5).Value = '  -0.03%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Cosmos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '11.39'
$ws.Cells.Item(26, 5).Value = '  +10.52%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '3.69'
$ws.Cells.Item(27, 5).Value = '  -2.62%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.41'
$ws.Cells.Item(28, 5).Value = '  -4.31%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -5.84%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '168.82'
$ws.Cells.Item(30, 5).Value = '  -1.54%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -8.80%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '20.13'
$ws.Cells.Item(32, 5).Value = '  -2.54%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.119'
$ws.Cells.Item(33, 5).Value = '  -2.04%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.68'
$ws.Cells.Item(34, 5).Value = '  +0.99%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0749'
$ws.Cells.Item(35, 5).Value = '  +3.87%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -3.26%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.56'
$ws.Cells.Item(37, 5).Value = '  -2.38%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.97'
$ws.Cells.Item(38, 5).Value = '  -0.56%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '24.62'
$ws.Cells.Item(39, 5).Value = '  -7.51%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0300'
$ws.Cells.Item(40, 5).Value = '  +4.41%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.17'
$ws.Cells.Item(41, 5).Value = '  -5.61%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.45'
$ws.Cells.Item(42, 5).Value = '  -9.50%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '11.44'
$ws.Cells.Item(43, 5).Value = '  -5.08%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '60.51'
$ws.Cells.Item(44, 5).Value = '  -12.70%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.77'
$ws.Cells.Item(45, 5).Value = '  -6.32%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.191'
$ws.Cells.Item(46, 5).Value = '  -8.49%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'BinanceUSD'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.00'
$ws.Cells.Item(47, 5).Value = '  +0.03%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0988'
$ws.Cells.Item(48, 5).Value = '  -2.61%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'FraxShare'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.40'
$ws.Cells.Item(49, 5).Value = '  -4.60%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.12'
$ws.Cells.Item(50, 5).Value = '  -2.20%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -4.11%  '
